## Actualización automática 2025-09-29 10:30:09
## Applies the sales updates for CASTRO ALCIVAR EDA MARIA across the three
## report sheets: VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

## --- Sheet "VENTAS POR GRUPO" --------------------------------------------
## SALAZAR BALLADARES MARIA ANGELICA now has INODOROS sales of 147.6
$wsGrupo.Range("H45").Value = 147.6

## TAPIA TAPIA ANGEL GUSTAVO now has SAL SOLUBLE sales of 869.53
$wsGrupo.Range("O53").Value = 869.53

## Row 57 "count of clients with sales" counters shift: one more client now
## has INODOROS sales (H57) and one more has SAL SOLUBLE sales (O57)
$wsGrupo.Range("H57").Value = "2 de 55"
$wsGrupo.Range("O57").Value = "1 de 55"

## --- Sheet "VENTA MENSUAL" -------------------------------------------------
## September (columna F) totals for the same two clients/groups
$wsMensual.Range("F45").Value = 4257.59
$wsMensual.Range("F53").Value = 869.53
$wsMensual.Range("F57").Value = 70559.19

## --- Sheet "CUMPLIMIENTO MENSUAL" ------------------------------------------
## INODOROS row (row 6): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D6").Value = 1605
$wsCumplimiento.Range("E6").Value = 721.0669451682102
$wsCumplimiento.Range("F6").Value = 0.6900059361292089

## SAL SOLUBLE row (row 14): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D14").Value = 869.53
$wsCumplimiento.Range("E14").Value = -379.726074704917
$wsCumplimiento.Range("F14").Value = 1.775261395621014

## TOTAL row (row 15): VENTA, POR CUMPLIR, CUMPLIMIENTO
$wsCumplimiento.Range("D15").Value = 89036.96000000001
$wsCumplimiento.Range("E15").Value = 9861.039925092748
$wsCumplimiento.Range("F15").Value = 0.9002908053493329
